$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.067569099508484
$ws.Range("D2").Value = 0.2972902393375034

$ws.Range("C3").Value = 0.9829279707501174
$ws.Range("D3").Value = 0.3363298452418406

$ws.Range("C4").Value = 0.8602695769791336
$ws.Range("D4").Value = 0.3989187246811412

$ws.Range("C5").Value = 0.6678845444242354
$ws.Range("D5").Value = 0.5111537394277872

$ws.Range("C6").Value = -0.1606875043405656
$ws.Range("D6").Value = 0.8738059162430534

$ws.Range("C7").Value = -0.1796496393343356
$ws.Range("D7").Value = 0.8590720150477924

$ws.Range("C8").Value = -0.3825527662245769
$ws.Range("D8").Value = 0.7057223204451546

$ws.Range("C9").Value = 0.01341334941452847
$ws.Range("D9").Value = 0.9894189146990597

$ws.Range("C10").Value = -0.1411662265816011
$ws.Range("D10").Value = 0.8890232749483302

$ws.Range("C11").Value = -0.1958786620863607
$ws.Range("D11").Value = 0.8465031228363822
